$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.317.31"
$ws.Range("E2").Value = "  +2.60%  "
$ws.Range("D3").Value = "2.063.65"
$ws.Range("E3").Value = "  +3.96%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.74"
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.613"
$ws.Range("E6").Value = "  +2.92%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.18"
$ws.Range("E7").Value = "  +7.25%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +3.59%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "59.19"
$ws.Range("E10").Value = "  +2.17%  "
$ws.Range("E11").Value = "  +1.95%  "
$ws.Range("E12").Value = "  +3.01%  "
$ws.Range("D13").Value = "2.369.02"
$ws.Range("E13").Value = "  +4.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.65"
$ws.Range("E14").Value = "  +3.67%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.98"
$ws.Range("E15").Value = "  +3.96%  "
$ws.Range("E16").Value = "  +3.37%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.18"
$ws.Range("E17").Value = "  +2.44%  "
$ws.Range("D18").Value = "2.086.36"
$ws.Range("E18").Value = "  +4.83%  "
$ws.Range("D19").Value = "37.576.56"
$ws.Range("E19").Value = "  +3.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.15"
$ws.Range("E20").Value = "  +17.25%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "68.97"
$ws.Range("E21").Value = "  +2.05%  "
$ws.Range("E22").Value = "  +1.56%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "226.40"
$ws.Range("E23").Value = "  +2.41%  "
$ws.Range("E24").Value = "  -0.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.44"
$ws.Range("E25").Value = "  +1.92%  "
$ws.Range("E26").Value = "  +0.96%  "
$ws.Range("B27").Value = "ImmutableX"
$ws.Range("C27").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.51"
$ws.Range("E27").Value = "  +14.53%  "
$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "164.73"
$ws.Range("E28").Value = "  +1.53%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.85"
$ws.Range("E29").Value = "  +2.44%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.19"
$ws.Range("E30").Value = "  +2.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.127"
$ws.Range("E31").Value = "  -1.34%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.118"
$ws.Range("E32").Value = "  +1.54%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.48"
$ws.Range("E33").Value = "  +2.78%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0622"
$ws.Range("E34").Value = "  +3.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.54"
$ws.Range("E35").Value = "  +9.23%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.50"
$ws.Range("E36").Value = "  +6.32%  "
$ws.Range("E37").Value = "  +1.83%  "
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("E39").Value = "  +0.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.85"
$ws.Range("E40").Value = "  +6.27%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0990"
$ws.Range("E41").Value = "  +6.29%  "
$ws.Range("E42").Value = "  -1.73%  "
$ws.Range("D43").Value = "1.472.69"
$ws.Range("E43").Value = "  +1.27%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "96.59"
$ws.Range("E44").Value = "  +8.46%  "
$ws.Range("B45").Value = "FTXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.29"
$ws.Range("E45").Value = "  +18.51%  "
$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.16"
$ws.Range("E46").Value = "  +5.80%  "
$ws.Range("E47").Value = "  +4.36%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "15.91"
$ws.Range("E48").Value = "  +6.48%  "
$ws.Range("E49").Value = "  +4.11%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.28"
$ws.Range("E50").Value = "  +6.92%  "
$ws.Range("E51").Value = "  +2.24%  "
